# Pemadanan Data.xlsx - edit replicating author's change:
#  - Kode_PKS value in A2 (sharedStrings "01732290") changed to "01732299"
#  - Active cell / selection on sheet moved from D4 to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Kode_PKS cell (A2). Use a leading apostrophe so the engine keeps
# treating/storing the value as text (preserves the existing quotePrefix
# cell style instead of re-classifying it as a plain number).
$ws.Range("A2").Value = "'01732299"

# Move the selection / active cell to D11, as in the authored edit.
$ws.Range("D11").Select()
